# Add two new rows of Python notes (row 12: if/else, row 13: min()/max())
# to the first worksheet, mirroring the existing Topic/Category/Details table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Text is carried as base64/UTF-8 so multi-line + unicode (…, 二/多元組, backticks)
# survive untouched through the PowerShell parser.
function Decode-B64 {
    param([string]$b64)
    $bytes = [System.Convert]::FromBase64String($b64)
    return [System.Text.Encoding]::UTF8.GetString($bytes)
}

$topicPython   = Decode-B64 "IyBQeXRob24="
$catIfElse     = Decode-B64 "IyMgaWYgZWxzZQ=="
$detailIfElse  = Decode-B64 "YGBgCmlmIGNvbmRpdGlvbjoKICAgIHN0YXRlbWVudCgpCmVsc2UgOgogICAgc3RhdGVtZW50KCk="
$catMinMax     = Decode-B64 "IyMgbWluKCkgJiBtYXgoKQ=="
$detailMinMax  = Decode-B64 "IyMjIGZpbmQgdGhlIG1pbi9tYXggaW4gYXJyYXkKYG1heChb4oCmXSlgCiMjIyBmaW5kIHRoZSBtaW4vbWF4IGluIGFycmF5IG9mIHR1cGxlcyDkuowv5aSa5YWD57WECmBtYXgoWyAoeDAsIHkwKSDigKYgKHhuLCB5bikgXSlgICAvL2luIHRoaXMgc2l0dWF0aW9uLCBtYXggYWNjb3JkaW5nIHRvIGZpcnN0IGVsZW1lbnQgaW4gdGhlIHR1cGxlcw=="

# Row 12: "# Python" / "## if else" / code block
$ws.Range("A12").Value2 = $topicPython
$ws.Range("B12").Value2 = $catIfElse
$ws.Range("C12").Value2 = $detailIfElse

# Row 13: "# Python" / "## min() & max()" / details
$ws.Range("A13").Value2 = $topicPython
$ws.Range("B13").Value2 = $catMinMax
$ws.Range("C13").Value2 = $detailMinMax

# Match styling of the other "Details"/"Category" rows (wrap text, column B/C look).
$ws.Range("B12:C13").WrapText = $true

# Row heights from the target layout.
$ws.Rows.Item(12).RowHeight = 63.75
$ws.Rows.Item(13).RowHeight = 51

# Update the view: scrolled down a couple more rows, new selection at C14.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C14").Select()
